$d = $word.ActiveDocument

function Insert-ParaXml($range, [string[]]$paragraphsXml) {
    $body = [string]::Join("`n", $paragraphsXml)
    $xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
$body
</w:body>
</w:document>
</pkg:xmlData></pkg:part></pkg:package>
"@
    $range.InsertXML($xml)
}

# --- Edit 5: "Planned work" heading paragraph (index 36) -----------------
# Split into the heading paragraph and a new paragraph with the planned-work
# text, moving the page-break run into the new paragraph.
$p36 = $d.Paragraphs.Item(36)
Insert-ParaXml $p36.Range @(
'<w:p><w:pPr><w:pStyle w:val="Heading3"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Planned work</w:t></w:r></w:p>',
'<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Make updates to our architecture, Finish the UML class diagram and the L2 documentation. Might also try to do one or two more UML state machine diagrams.</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:br w:type="page"/></w:r></w:p>'
)

# --- Edit 4: "Major design decisions" heading paragraph (index 32) -------
# Add a new paragraph right after it noting nothing changed this week.
$p32 = $d.Paragraphs.Item(32)
Insert-ParaXml $p32.Range @(
'<w:p><w:pPr><w:pStyle w:val="Heading3"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Major design decisions</w:t></w:r></w:p>',
'<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>No major design decisions have been made this week.</w:t></w:r></w:p>'
)

# --- Edit 3: Work summary paragraph (index 31) ----------------------------
# Merge the split "Pacman" runs back into one run and append a new sentence
# about starting the L2 documentation work.
$p31 = $d.Paragraphs.Item(31)
Insert-ParaXml $p31.Range @(
'<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>We reworked our time chart, making it more detailed. We continued to construct more smaller state machine diagrams, one for the different states of Pacman and one for the different states of the screen handler.</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> We also started programming the skeleton for the project.</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Started work on the Documentation for the L2 delivery, making small description texts for each class in our project.</w:t></w:r></w:p>'
)

# --- Edit 2: "Calle Ketola" footer paragraph (index 28) -------------------
$p28 = $d.Paragraphs.Item(28)
Insert-ParaXml $p28.Range @(
'<w:p><w:pPr><w:pStyle w:val="Footer"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Calle Ketola - cake10@student.bth.se</w:t></w:r></w:p>'
)

# --- Edit 1: "Rasmus Tilljander" footer paragraph (index 26) --------------
$p26 = $d.Paragraphs.Item(26)
Insert-ParaXml $p26.Range @(
'<w:p><w:pPr><w:pStyle w:val="Footer"/></w:pPr><w:r><w:t>Rasmus Tilljander - rati10@student.bth.se</w:t></w:r></w:p>'
)

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
